$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# (values like "5.50" or "0.0240" would otherwise be coerced to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.553.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.477.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.21"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.42%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000184"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.909.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.346.86"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.487.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.91"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +11.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "639.08"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +15.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.85"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000105"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +13.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.637.15"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.52"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +10.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.51"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.991"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.23"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.387"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.50"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.88"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +18.03%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "151.67"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.80"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0553"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.611"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0929"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.22%  "

# Row 51: ONDO -> WhiteBITCoin (full row replacement)
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.61"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.51%  "
